# Auto-generated: update odds values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16, 8).Value = 3.4
$ws.Cells.Item(16, 10).Value = 2.5
$ws.Cells.Item(16, 11).Value = 2.05
$ws.Cells.Item(16, 13).Value = 1.08
$ws.Cells.Item(16, 14).Value = 8
$ws.Cells.Item(16, 15).Value = 1.4
$ws.Cells.Item(16, 16).Value = 2.75
$ws.Cells.Item(16, 17).Value = 2.25
$ws.Cells.Item(16, 18).Value = 1.62
$ws.Cells.Item(16, 19).Value = 1.5
$ws.Cells.Item(16, 20).Value = 2.5
$ws.Cells.Item(16, 23).Value = 6
$ws.Cells.Item(16, 24).Value = 8
$ws.Cells.Item(16, 27).Value = 17
$ws.Cells.Item(16, 29).Value = 7.5
$ws.Cells.Item(16, 36).Value = 15
$ws.Cells.Item(16, 41).Value = 10
$ws.Cells.Item(16, 42).Value = 23
$ws.Cells.Item(16, 43).Value = 34
$ws.Cells.Item(16, 46).Value = 2.5
$ws.Cells.Item(16, 53).Value = 34
$ws.Cells.Item(16, 55).Value = 126
$ws.Cells.Item(16, 56).Value = 301

# Row 17
$ws.Cells.Item(17, 21).Value = 2.2
$ws.Cells.Item(17, 22).Value = 1.62
$ws.Cells.Item(17, 25).Value = 9.5
$ws.Cells.Item(17, 29).Value = 7
$ws.Cells.Item(17, 31).Value = 21
$ws.Cells.Item(17, 32).Value = 81
$ws.Cells.Item(17, 35).Value = 21
$ws.Cells.Item(17, 45).Value = 251
$ws.Cells.Item(17, 55).Value = 151

# Row 18
$ws.Cells.Item(18, 7).Value = 2.88
$ws.Cells.Item(18, 8).Value = 2.88
$ws.Cells.Item(18, 9).Value = 2.7
$ws.Cells.Item(18, 10).Value = 3.75
$ws.Cells.Item(18, 12).Value = 3.75
$ws.Cells.Item(18, 23).Value = 6
$ws.Cells.Item(18, 24).Value = 12
$ws.Cells.Item(18, 26).Value = 29
$ws.Cells.Item(18, 27).Value = 34
$ws.Cells.Item(18, 34).Value = 6
$ws.Cells.Item(18, 35).Value = 11
$ws.Cells.Item(18, 36).Value = 12
$ws.Cells.Item(18, 37).Value = 29
$ws.Cells.Item(18, 38).Value = 29
$ws.Cells.Item(18, 40).Value = 4.75
$ws.Cells.Item(18, 41).Value = 19
$ws.Cells.Item(18, 42).Value = 41
$ws.Cells.Item(18, 43).Value = 67
$ws.Cells.Item(18, 44).Value = 126
$ws.Cells.Item(18, 51).Value = 4.5
$ws.Cells.Item(18, 52).Value = 19
$ws.Cells.Item(18, 53).Value = 34

# Row 19
$ws.Cells.Item(19, 7).Value = 1.85
$ws.Cells.Item(19, 8).Value = 3.25
$ws.Cells.Item(19, 9).Value = 4.75
$ws.Cells.Item(19, 10).Value = 2.6
$ws.Cells.Item(19, 12).Value = 6
$ws.Cells.Item(19, 13).Value = 1.13
$ws.Cells.Item(19, 14).Value = 6
$ws.Cells.Item(19, 17).Value = 2.88
$ws.Cells.Item(19, 18).Value = 1.4
$ws.Cells.Item(19, 21).Value = 2.5
$ws.Cells.Item(19, 22).Value = 1.5
$ws.Cells.Item(19, 24).Value = 7
$ws.Cells.Item(19, 27).Value = 21
$ws.Cells.Item(19, 34).Value = 9
$ws.Cells.Item(19, 35).Value = 23
$ws.Cells.Item(19, 37).Value = 51
$ws.Cells.Item(19, 40).Value = 3.5
$ws.Cells.Item(19, 43).Value = 41
$ws.Cells.Item(19, 44).Value = 81
$ws.Cells.Item(19, 47).Value = 10
$ws.Cells.Item(19, 48).Value = 81
$ws.Cells.Item(19, 51).Value = 6.5
$ws.Cells.Item(19, 55).Value = 201

# Row 23
$ws.Cells.Item(23, 7).Value = 3.6
$ws.Cells.Item(23, 8).Value = 2.9
$ws.Cells.Item(23, 9).Value = 2.3
$ws.Cells.Item(23, 10).Value = 4.33
$ws.Cells.Item(23, 12).Value = 3.2
$ws.Cells.Item(23, 13).Value = 1.11
$ws.Cells.Item(23, 14).Value = 6.5
$ws.Cells.Item(23, 19).Value = 1.62
$ws.Cells.Item(23, 20).Value = 2.2
$ws.Cells.Item(23, 22).Value = 1.58
$ws.Cells.Item(23, 23).Value = 8
$ws.Cells.Item(23, 26).Value = 41
$ws.Cells.Item(23, 28).Value = 51
$ws.Cells.Item(23, 32).Value = 81
$ws.Cells.Item(23, 35).Value = 9.5
$ws.Cells.Item(23, 36).Value = 10
$ws.Cells.Item(23, 37).Value = 21
$ws.Cells.Item(23, 43).Value = 81
$ws.Cells.Item(23, 44).Value = 126
$ws.Cells.Item(23, 46).Value = 2.2
$ws.Cells.Item(23, 47).Value = 9.5
$ws.Cells.Item(23, 51).Value = 4

# Row 29
$ws.Cells.Item(29, 7).Value = 1.57
$ws.Cells.Item(29, 9).Value = 6
$ws.Cells.Item(29, 17).Value = 2.1
$ws.Cells.Item(29, 18).Value = 1.7
$ws.Cells.Item(29, 29).Value = 8
$ws.Cells.Item(29, 30).Value = 7
$ws.Cells.Item(29, 31).Value = 21
$ws.Cells.Item(29, 34).Value = 13

# Row 45
$ws.Cells.Item(45, 7).Value = 1.48
$ws.Cells.Item(45, 8).Value = 4.5
$ws.Cells.Item(45, 9).Value = 7
$ws.Cells.Item(45, 10).Value = 2
$ws.Cells.Item(45, 11).Value = 2.4
$ws.Cells.Item(45, 13).Value = 1.04
$ws.Cells.Item(45, 14).Value = 13
$ws.Cells.Item(45, 17).Value = 1.67
$ws.Cells.Item(45, 18).Value = 2.15
$ws.Cells.Item(45, 19).Value = 1.3
$ws.Cells.Item(45, 20).Value = 3.4
$ws.Cells.Item(45, 21).Value = 1.8
$ws.Cells.Item(45, 22).Value = 1.95
$ws.Cells.Item(45, 25).Value = 8
$ws.Cells.Item(45, 26).Value = 10
$ws.Cells.Item(45, 29).Value = 13
$ws.Cells.Item(45, 30).Value = 8.5
$ws.Cells.Item(45, 34).Value = 19
$ws.Cells.Item(45, 36).Value = 21
$ws.Cells.Item(45, 41).Value = 7
$ws.Cells.Item(45, 46).Value = 3.4
$ws.Cells.Item(45, 56).Value = 201

# Row 46
$ws.Cells.Item(46, 7).Value = 1.52
$ws.Cells.Item(46, 8).Value = 4
$ws.Cells.Item(46, 9).Value = 5.5
$ws.Cells.Item(46, 10).Value = 2.05
$ws.Cells.Item(46, 11).Value = 2.25
$ws.Cells.Item(46, 12).Value = 5.4
$ws.Cells.Item(46, 13).Value = 1.01
$ws.Cells.Item(46, 14).Value = 11
$ws.Cells.Item(46, 19).Value = 1.34
$ws.Cells.Item(46, 20).Value = 3.1
$ws.Cells.Item(46, 21).Value = 1.83
$ws.Cells.Item(46, 22).Value = 1.78
$ws.Cells.Item(46, 23).Value = 6.8
$ws.Cells.Item(46, 24).Value = 7
$ws.Cells.Item(46, 26).Value = 10.5
$ws.Cells.Item(46, 28).Value = 27
$ws.Cells.Item(46, 30).Value = 7.9
$ws.Cells.Item(46, 31).Value = 17.5
$ws.Cells.Item(46, 32).Value = 80
$ws.Cells.Item(46, 33).Value = 700
$ws.Cells.Item(46, 34).Value = 15
$ws.Cells.Item(46, 35).Value = 35
$ws.Cells.Item(46, 36).Value = 17.5
$ws.Cells.Item(46, 37).Value = 110
$ws.Cells.Item(46, 38).Value = 55
$ws.Cells.Item(46, 39).Value = 55
$ws.Cells.Item(46, 40).Value = 3.3
$ws.Cells.Item(46, 41).Value = 7.1
$ws.Cells.Item(46, 42).Value = 17
$ws.Cells.Item(46, 43).Value = 22
$ws.Cells.Item(46, 44).Value = 55
$ws.Cells.Item(46, 47).Value = 7.9
$ws.Cells.Item(46, 48).Value = 75
$ws.Cells.Item(46, 51).Value = 6.9
$ws.Cells.Item(46, 52).Value = 30
$ws.Cells.Item(46, 53).Value = 35
$ws.Cells.Item(46, 54).Value = 200

# Row 47
$ws.Cells.Item(47, 8).Value = 4.6
$ws.Cells.Item(47, 9).Value = 7
$ws.Cells.Item(47, 11).Value = 2.42
$ws.Cells.Item(47, 12).Value = 6.2
$ws.Cells.Item(47, 17).Value = 1.55
$ws.Cells.Item(47, 18).Value = 2.15
$ws.Cells.Item(47, 22).Value = 1.82
$ws.Cells.Item(47, 23).Value = 7.8
$ws.Cells.Item(47, 24).Value = 7
$ws.Cells.Item(47, 29).Value = 14.5
$ws.Cells.Item(47, 33).Value = 600
$ws.Cells.Item(47, 34).Value = 21
$ws.Cells.Item(47, 37).Value = 150
$ws.Cells.Item(47, 39).Value = 60
$ws.Cells.Item(47, 42).Value = 15.5
$ws.Cells.Item(47, 45).Value = 200
$ws.Cells.Item(47, 47).Value = 8
$ws.Cells.Item(47, 48).Value = 70
$ws.Cells.Item(47, 51).Value = 8
$ws.Cells.Item(47, 52).Value = 37

# Row 61
$ws.Cells.Item(61, 15).Value = 1.36
$ws.Cells.Item(61, 16).Value = 3
$ws.Cells.Item(61, 17).Value = 2.2
$ws.Cells.Item(61, 18).Value = 1.65
$ws.Cells.Item(61, 33).Value = 1000

# Row 86
$ws.Cells.Item(86, 14).Value = 7.6
$ws.Cells.Item(86, 18).Value = 1.88
$ws.Cells.Item(86, 19).Value = 1.39
$ws.Cells.Item(86, 20).Value = 2.77
$ws.Cells.Item(86, 22).Value = 1.85
$ws.Cells.Item(86, 23).Value = 6.9
$ws.Cells.Item(86, 29).Value = 7.6
$ws.Cells.Item(86, 43).Value = 27
$ws.Cells.Item(86, 44).Value = 60
$ws.Cells.Item(86, 46).Value = 2.77
$ws.Cells.Item(86, 47).Value = 7.8
$ws.Cells.Item(86, 48).Value = 75

# Row 88
$ws.Cells.Item(88, 7).Value = 3.2
$ws.Cells.Item(88, 8).Value = 2.95
$ws.Cells.Item(88, 10).Value = 3.85
$ws.Cells.Item(88, 11).Value = 1.88
$ws.Cells.Item(88, 12).Value = 3.05
$ws.Cells.Item(88, 14).Value = 6.9
$ws.Cells.Item(88, 23).Value = 8.5
$ws.Cells.Item(88, 24).Value = 16.5
$ws.Cells.Item(88, 25).Value = 11.25
$ws.Cells.Item(88, 27).Value = 30
$ws.Cells.Item(88, 28).Value = 40
$ws.Cells.Item(88, 29).Value = 7.4
$ws.Cells.Item(88, 30).Value = 5.8
$ws.Cells.Item(88, 34).Value = 6.4
$ws.Cells.Item(88, 35).Value = 10.25
$ws.Cells.Item(88, 36).Value = 9.25
$ws.Cells.Item(88, 37).Value = 23
$ws.Cells.Item(88, 38).Value = 21
$ws.Cells.Item(88, 39).Value = 35
$ws.Cells.Item(88, 40).Value = 4.9
$ws.Cells.Item(88, 41).Value = 19
$ws.Cells.Item(88, 42).Value = 29
$ws.Cells.Item(88, 43).Value = 100
$ws.Cells.Item(88, 44).Value = 150
$ws.Cells.Item(88, 45).Value = 450
$ws.Cells.Item(88, 52).Value = 13
$ws.Cells.Item(88, 53).Value = 25
$ws.Cells.Item(88, 54).Value = 60
$ws.Cells.Item(88, 55).Value = 120
$ws.Cells.Item(88, 56).Value = 400
